# Restore the "Artificial Intelligence and Machine Learning" multi-industry
# template content in place of the "Information Technology" variant, and
# drop the leftover truly-blank spacer rows that shouldn't be persisted.
#
# This mirrors commit 168d9c4 ("RESTORE: Recover all 973 original
# multi-industry template files") for the IT_Comprehensive_Budget workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Instructions & User Guide"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning Comprehensive Budget - User Guide & Instructions"
$ws.Cells.Item(56, 1).Value = "📋 ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING PROJECT OVERVIEW"
$ws.Cells.Item(59, 2).Value = "Data Scientists, ML Engineers, AI Architects, DevOps Engineers..."

# ---------------------------------------------------------------------
# Sheet 2: "Budget Summary"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Executive Budget Summary"

# ---------------------------------------------------------------------
# Sheet 3: "Resources"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Resources Budget"
$ws.Cells.Item(4, 1).Value = "Data Scientists"
$ws.Cells.Item(5, 1).Value = "ML Engineers"
$ws.Cells.Item(9, 1).Value = "Business Analysts"

# ---------------------------------------------------------------------
# Sheet 4: "Logistics"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Logistics Budget"

# ---------------------------------------------------------------------
# Sheet 5: "Technology"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Technology Budget"

# ---------------------------------------------------------------------
# Sheet 6: "Training"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Training Budget"
$ws.Cells.Item(4, 1).Value = "AI/ML Certification Programs"

# ---------------------------------------------------------------------
# Sheet 7: "Contingency"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Contingency Budget"

# ---------------------------------------------------------------------
# Sheet 8: "Timeline"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")

$ws.Cells.Item(1, 1).Value = "Artificial Intelligence and Machine Learning - Budget Timeline"
